$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column L, duplicating the 2020 figures already present in
# column K (value + formatting) for the header/year row and the data row.
$ws.Range("K3").Copy($ws.Range("L3"))
$ws.Range("K4").Copy($ws.Range("L4"))

# Match the author's final selection.
$ws.Range("L10").Select()
